# Update the "想去人数" (want-to-go count) column F values on the
# "展览", "演出" and "全部类型" sheets to match the newly generated data.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1106
$ws1.Range("F5").Value  = 992
$ws1.Range("F7").Value  = 547
$ws1.Range("F9").Value  = 53
$ws1.Range("F15").Value = 662
$ws1.Range("F16").Value = 151
$ws1.Range("F18").Value = 25
$ws1.Range("F20").Value = 327
$ws1.Range("F28").Value = 308

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value  = 248

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1106
$ws4.Range("F6").Value  = 992
$ws4.Range("F8").Value  = 547
$ws4.Range("F10").Value = 53
$ws4.Range("F17").Value = 662
$ws4.Range("F18").Value = 151
$ws4.Range("F21").Value = 25
$ws4.Range("F25").Value = 327
$ws4.Range("F27").Value = 248
$ws4.Range("F28").Value = 248
$ws4.Range("F36").Value = 308
